$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap row pairs: F:V content exchanged, A:E left untouched ---
$rowA = $ws.Range("F42:V42").Value2
$rowB = $ws.Range("F43:V43").Value2
for ($i = 1; $i -le 17; $i++) {
    $ws.Cells.Item(42, 5+$i).Value = $rowB[1,$i]
    $ws.Cells.Item(43, 5+$i).Value = $rowA[1,$i]
}

$rowA = $ws.Range("F44:V44").Value2
$rowB = $ws.Range("F45:V45").Value2
for ($i = 1; $i -le 17; $i++) {
    $ws.Cells.Item(44, 5+$i).Value = $rowB[1,$i]
    $ws.Cells.Item(45, 5+$i).Value = $rowA[1,$i]
}

$rowA = $ws.Range("F84:V84").Value2
$rowB = $ws.Range("F85:V85").Value2
for ($i = 1; $i -le 17; $i++) {
    $ws.Cells.Item(84, 5+$i).Value = $rowB[1,$i]
    $ws.Cells.Item(85, 5+$i).Value = $rowA[1,$i]
}

$rowA = $ws.Range("F96:V96").Value2
$rowB = $ws.Range("F97:V97").Value2
for ($i = 1; $i -le 17; $i++) {
    $ws.Cells.Item(96, 5+$i).Value = $rowB[1,$i]
    $ws.Cells.Item(97, 5+$i).Value = $rowA[1,$i]
}

# --- Append new rows 111-117 (matches onward) ---
$ws.Range("A110:V110").Copy()
$ws.Range("A111:V111").PasteSpecial(-4122)
$ws.Cells.Item(111, 1).Value = 110
$ws.Cells.Item(111, 2).Value = "turkey"
$ws.Cells.Item(111, 3).Value = "1-lig"
$ws.Cells.Item(111, 4).Value = "2023-2024"
$ws.Cells.Item(111, 5).Value = 45255.47916666666
$ws.Cells.Item(111, 6).Value = "Bandirmaspor"
$ws.Cells.Item(111, 7).Value = 0
$ws.Cells.Item(111, 8).Value = "Sanliurfaspor"
$ws.Cells.Item(111, 9).Value = 1
$ws.Cells.Item(111, 10).Value = 1.58
$ws.Cells.Item(111, 11).Value = "17/11/2023 11:42"
$ws.Cells.Item(111, 12).Value = 1.3
$ws.Cells.Item(111, 13).Value = "25/11/2023 11:13"
$ws.Cells.Item(111, 14).Value = 4.05
$ws.Cells.Item(111, 15).Value = "17/11/2023 11:42"
$ws.Cells.Item(111, 16).Value = 5.49
$ws.Cells.Item(111, 17).Value = "25/11/2023 11:20"
$ws.Cells.Item(111, 18).Value = 5.6
$ws.Cells.Item(111, 19).Value = "17/11/2023 11:42"
$ws.Cells.Item(111, 20).Value = 9.88
$ws.Cells.Item(111, 21).Value = "25/11/2023 11:20"
$ws.Cells.Item(111, 22).Value = "https://www.betexplorer.com/football/turkey/1-lig/bandirmaspor-sanliurfaspor/d6jXgAs8/"

$ws.Range("A111:V111").Copy()
$ws.Range("A112:V112").PasteSpecial(-4122)
$ws.Cells.Item(112, 1).Value = 111
$ws.Cells.Item(112, 2).Value = "turkey"
$ws.Cells.Item(112, 3).Value = "1-lig"
$ws.Cells.Item(112, 4).Value = "2023-2024"
$ws.Cells.Item(112, 5).Value = 45255.47916666666
$ws.Cells.Item(112, 6).Value = "Bodrumspor"
$ws.Cells.Item(112, 7).Value = 0
$ws.Cells.Item(112, 8).Value = "Adanaspor AS"
$ws.Cells.Item(112, 9).Value = 1
$ws.Cells.Item(112, 10).Value = 1.47
$ws.Cells.Item(112, 11).Value = "17/11/2023 11:42"
$ws.Cells.Item(112, 12).Value = 1.43
$ws.Cells.Item(112, 13).Value = "25/11/2023 11:25"
$ws.Cells.Item(112, 14).Value = 4.42
$ws.Cells.Item(112, 15).Value = "17/11/2023 11:42"
$ws.Cells.Item(112, 16).Value = 4.53
$ws.Cells.Item(112, 17).Value = "25/11/2023 11:29"
$ws.Cells.Item(112, 18).Value = 6.43
$ws.Cells.Item(112, 19).Value = "17/11/2023 11:42"
$ws.Cells.Item(112, 20).Value = 7.5
$ws.Cells.Item(112, 21).Value = "25/11/2023 11:29"
$ws.Cells.Item(112, 22).Value = "https://www.betexplorer.com/football/turkey/1-lig/bodrumspor-adanaspor-as/EsRnnbO1/"

$ws.Range("A112:V112").Copy()
$ws.Range("A113:V113").PasteSpecial(-4122)
$ws.Cells.Item(113, 1).Value = 112
$ws.Cells.Item(113, 2).Value = "turkey"
$ws.Cells.Item(113, 3).Value = "1-lig"
$ws.Cells.Item(113, 4).Value = "2023-2024"
$ws.Cells.Item(113, 5).Value = 45255.58333333334
$ws.Cells.Item(113, 6).Value = "Umraniyespor"
$ws.Cells.Item(113, 7).Value = 1
$ws.Cells.Item(113, 8).Value = "Boluspor"
$ws.Cells.Item(113, 9).Value = 2
$ws.Cells.Item(113, 10).Value = 1.99
$ws.Cells.Item(113, 11).Value = "17/11/2023 14:42"
$ws.Cells.Item(113, 12).Value = 2.18
$ws.Cells.Item(113, 13).Value = "25/11/2023 13:34"
$ws.Cells.Item(113, 14).Value = 3.44
$ws.Cells.Item(113, 15).Value = "17/11/2023 14:42"
$ws.Cells.Item(113, 16).Value = 3.2
$ws.Cells.Item(113, 17).Value = "25/11/2023 13:34"
$ws.Cells.Item(113, 18).Value = 3.81
$ws.Cells.Item(113, 19).Value = "17/11/2023 14:42"
$ws.Cells.Item(113, 20).Value = 3.64
$ws.Cells.Item(113, 21).Value = "25/11/2023 13:34"
$ws.Cells.Item(113, 22).Value = "https://www.betexplorer.com/football/turkey/1-lig/umraniyespor-boluspor/p8Mrmv9e/"

$ws.Range("A113:V113").Copy()
$ws.Range("A114:V114").PasteSpecial(-4122)
$ws.Cells.Item(114, 1).Value = 113
$ws.Cells.Item(114, 2).Value = "turkey"
$ws.Cells.Item(114, 3).Value = "1-lig"
$ws.Cells.Item(114, 4).Value = "2023-2024"
$ws.Cells.Item(114, 5).Value = 45255.70833333334
$ws.Cells.Item(114, 6).Value = "Giresunspor"
$ws.Cells.Item(114, 7).Value = 1
$ws.Cells.Item(114, 8).Value = "Kocaelispor"
$ws.Cells.Item(114, 9).Value = 4
$ws.Cells.Item(114, 10).Value = 4.16
$ws.Cells.Item(114, 11).Value = "17/11/2023 17:43"
$ws.Cells.Item(114, 12).Value = 5.71
$ws.Cells.Item(114, 13).Value = "25/11/2023 16:51"
$ws.Cells.Item(114, 14).Value = 3.64
$ws.Cells.Item(114, 15).Value = "17/11/2023 17:43"
$ws.Cells.Item(114, 16).Value = 3.97
$ws.Cells.Item(114, 17).Value = "25/11/2023 16:51"
$ws.Cells.Item(114, 18).Value = 1.85
$ws.Cells.Item(114, 19).Value = "17/11/2023 17:43"
$ws.Cells.Item(114, 20).Value = 1.61
$ws.Cells.Item(114, 21).Value = "25/11/2023 16:51"
$ws.Cells.Item(114, 22).Value = "https://www.betexplorer.com/football/turkey/1-lig/giresunspor-kocaelispor/rNhqi8RQ/"

$ws.Range("A114:V114").Copy()
$ws.Range("A115:V115").PasteSpecial(-4122)
$ws.Cells.Item(115, 1).Value = 114
$ws.Cells.Item(115, 2).Value = "turkey"
$ws.Cells.Item(115, 3).Value = "1-lig"
$ws.Cells.Item(115, 4).Value = "2023-2024"
$ws.Cells.Item(115, 5).Value = 45256.47916666666
$ws.Cells.Item(115, 6).Value = "Corum"
$ws.Cells.Item(115, 7).Value = 3
$ws.Cells.Item(115, 8).Value = "Altay"
$ws.Cells.Item(115, 9).Value = 0
$ws.Cells.Item(115, 10).Value = 1.34
$ws.Cells.Item(115, 11).Value = "18/11/2023 11:42"
$ws.Cells.Item(115, 12).Value = 1.38
$ws.Cells.Item(115, 13).Value = "26/11/2023 11:20"
$ws.Cells.Item(115, 14).Value = 4.99
$ws.Cells.Item(115, 15).Value = "18/11/2023 11:42"
$ws.Cells.Item(115, 16).Value = 5.35
$ws.Cells.Item(115, 17).Value = "26/11/2023 11:28"
$ws.Cells.Item(115, 18).Value = 8.66
$ws.Cells.Item(115, 19).Value = "18/11/2023 11:42"
$ws.Cells.Item(115, 20).Value = 7.13
$ws.Cells.Item(115, 21).Value = "26/11/2023 11:26"
$ws.Cells.Item(115, 22).Value = "https://www.betexplorer.com/football/turkey/1-lig/corum-fk-altay/hYPjoIw8/"

$ws.Range("A115:V115").Copy()
$ws.Range("A116:V116").PasteSpecial(-4122)
$ws.Cells.Item(116, 1).Value = 115
$ws.Cells.Item(116, 2).Value = "turkey"
$ws.Cells.Item(116, 3).Value = "1-lig"
$ws.Cells.Item(116, 4).Value = "2023-2024"
$ws.Cells.Item(116, 5).Value = 45256.58333333334
$ws.Cells.Item(116, 6).Value = "Eyupspor"
$ws.Cells.Item(116, 7).Value = 2
$ws.Cells.Item(116, 8).Value = "Manisa FK"
$ws.Cells.Item(116, 9).Value = 1
$ws.Cells.Item(116, 10).Value = 1.51
$ws.Cells.Item(116, 11).Value = "18/11/2023 14:42"
$ws.Cells.Item(116, 12).Value = 1.56
$ws.Cells.Item(116, 13).Value = "26/11/2023 13:58"
$ws.Cells.Item(116, 14).Value = 4.23
$ws.Cells.Item(116, 15).Value = "18/11/2023 14:42"
$ws.Cells.Item(116, 16).Value = 4.13
$ws.Cells.Item(116, 17).Value = "26/11/2023 13:58"
$ws.Cells.Item(116, 18).Value = 6.2
$ws.Cells.Item(116, 19).Value = "18/11/2023 14:42"
$ws.Cells.Item(116, 20).Value = 6.04
$ws.Cells.Item(116, 21).Value = "26/11/2023 13:58"
$ws.Cells.Item(116, 22).Value = "https://www.betexplorer.com/football/turkey/1-lig/eyupspor-manisa-fk/A9fygUdE/"

$ws.Range("A116:V116").Copy()
$ws.Range("A117:V117").PasteSpecial(-4122)
$ws.Cells.Item(117, 1).Value = 116
$ws.Cells.Item(117, 2).Value = "turkey"
$ws.Cells.Item(117, 3).Value = "1-lig"
$ws.Cells.Item(117, 4).Value = "2023-2024"
$ws.Cells.Item(117, 5).Value = 45256.70833333334
$ws.Cells.Item(117, 6).Value = "Goztepe"
$ws.Cells.Item(117, 7).Value = 2
$ws.Cells.Item(117, 8).Value = "Tuzlaspor"
$ws.Cells.Item(117, 9).Value = 0
$ws.Cells.Item(117, 10).Value = 1.5
$ws.Cells.Item(117, 11).Value = "18/11/2023 18:13"
$ws.Cells.Item(117, 12).Value = 1.51
$ws.Cells.Item(117, 13).Value = "26/11/2023 16:58"
$ws.Cells.Item(117, 14).Value = 4.22
$ws.Cells.Item(117, 15).Value = "18/11/2023 18:13"
$ws.Cells.Item(117, 16).Value = 4.08
$ws.Cells.Item(117, 17).Value = "26/11/2023 16:58"
$ws.Cells.Item(117, 18).Value = 6.42
$ws.Cells.Item(117, 19).Value = "18/11/2023 18:13"
$ws.Cells.Item(117, 20).Value = 6.93
$ws.Cells.Item(117, 21).Value = "26/11/2023 16:58"
$ws.Cells.Item(117, 22).Value = "https://www.betexplorer.com/football/turkey/1-lig/goztepe-tuzlaspor/4WguhlCK/"

$excel.CutCopyMode = 0